$d = $word.ActiveDocument

# 1) Add the missing office number "8487" into the Ofício reference.
$d.Content.Find.Execute(" , datado de 04/01/2024, oriundo da ", $true, $false, $false, $false, $false, `
    $true, 1, $false, " 8487, datado de 04/01/2024, oriundo da ", 2)

# 2) Change the motorcycle color from Azul to Verde.
$d.Content.Find.Execute(", ano de fabricação/modelo 2010/789 de cor Azul", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", ano de fabricação/modelo 2010/789 de cor Verde", 2)

# 3) Fix the typo PERICIDA -> PERICIADA in the caption under the photos.
$d.Content.Find.Execute("MOTOCICLETA PERICIDA", $true, $false, $false, $false, $false, `
    $true, 1, $false, "MOTOCICLETA PERICIADA", 2)

# 4) Append the new "DO EXAME" section at the end of the document (after
#    the "MOTOCICLETA PERICIADA" caption paragraph, before the sectPr).
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newParagraphsXml = @"
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="justify"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:eastAsia="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:b w:val="1"/>
      <w:bCs w:val="1"/>
    </w:rPr>
    <w:t xml:space="preserve">DO EXAME</w:t>
  </w:r>
</w:p>
<w:p $ns/>
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="justify"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:eastAsia="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
    </w:rPr>
    <w:t xml:space="preserve">Com relação às numerações identificadoras da motocicletaforam observados:</w:t>
  </w:r>
</w:p>
<w:p $ns/>
<w:p $ns>
  <w:pPr>
    <w:pStyle w:val="justify"/>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:eastAsia="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:b w:val="0"/>
      <w:bCs w:val="0"/>
    </w:rPr>
    <w:t xml:space="preserve">a) número do chassi: esta numeração na xxxxxxxxx periciada se encontra gravada no xxxxxxxxxxxxxxxxx. Ao exame de referido suporte, após a devida limpeza, foi verificada a gravação da sequência alfanumérica xxxxxxxxxxxxxxx, a qual apresenta-se íntegra, sem sinais ou vestígios de adulteração.</w:t>
  </w:r>
</w:p>
<w:p $ns/>
<w:p $ns>
  <w:pPr/>
  <w:r>
    <w:rPr/>
    <w:t xml:space="preserve"></w:t>
  </w:r>
</w:p>
"@

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($newParagraphsXml)
